$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Status" header (C5) to "Completed"
$ws.Range("C5").Value = "Completed"

# Mark letter "q" (row 22) as completed: set Status to "Yes" and copy the
# same fill formatting used by the other "Yes" rows (e.g. C6)
$ws.Range("C22").Value = "Yes"
$ws.Range("C6").Copy()
$ws.Range("C22").PasteSpecial(-4122)

# Record the completion date for row 22, matching the date formatting
# used elsewhere in column D (e.g. D6)
$ws.Range("D6").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 43613

# Clear clipboard marching ants / selection artifacts
$excel.CutCopyMode = $false

# Update the active selection to C5, as recorded in the saved view state
$ws.Range("C5").Select()
